$wb = $excel.ActiveWorkbook

# --- "Chart" sheet: append 3 new daily rows (2025-12-02 .. 2025-12-04) ---
$chart = $wb.Worksheets.Item("Chart")

# Row 60: 2025-12-02
# Leading apostrophe forces the date-shaped string to stay plain text
# (matching how every other "Date" column cell is stored as a shared string).
$chart.Range("A60").Value = "'2025-12-02"
$chart.Range("B60").Value = 24
$chart.Range("C60").Value = 1
$chart.Range("D60").Value = 0

# Row 61: 2025-12-03
$chart.Range("A61").Value = "'2025-12-03"
$chart.Range("B61").Value = 24
$chart.Range("C61").Value = 1
$chart.Range("D61").Value = 0

# Row 62: 2025-12-04
$chart.Range("A62").Value = "'2025-12-04"
$chart.Range("B62").Value = 24
$chart.Range("C62").Value = 1
$chart.Range("D62").Value = 0

# --- "Table" sheet: refreshed validation-failure tally (23 -> 24) ---
$table = $wb.Worksheets.Item("Table")
$table.Range("C2").Value = 24
